$d = $word.ActiveDocument

# Locate the "GPA:  3.7" run and capture its range.
$rng = $d.Content
$found = $rng.Find.Execute("GPA:  3.7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $gpaEnd = $rng.End

    # Replace the trailing "7" with "6".
    $lastChar = $d.Range($gpaEnd - 1, $gpaEnd)
    $lastChar.Text = "6"

    # Force the newly written character into its own run (distinct from the
    # "GPA:  3." run that precedes it) by toggling a character formatting
    # property off and back on; this mirrors how Word keeps separately
    # authored runs from silently re-merging even when their resulting
    # properties are identical.
    $splitRng = $d.Range($gpaEnd - 1, $gpaEnd)
    $splitRng.Font.Bold = $false
    $splitRng2 = $d.Range($gpaEnd - 1, $gpaEnd)
    $splitRng2.Font.Bold = $true
}
